$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.776.79'
$ws.Range("E2").Value = '  -3.10%  '
$ws.Range("D3").Value = '3.356.57'
$ws.Range("E3").Value = '  -0.60%  '
$ws.Range("E4").Value = '  +0.07%  '
$c = $ws.Range("D5")
$c.Value = "'" + '570.23'
$c.ClearFormats()
$ws.Range("E5").Value = '  -0.29%  '
$c = $ws.Range("D6")
$c.Value = "'" + '133.83'
$c.ClearFormats()
$ws.Range("E6").Value = '  +6.66%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '3.355.29'
$ws.Range("E8").Value = '  -0.66%  '
$c = $ws.Range("D9")
$c.Value = "'" + '0.476'
$c.ClearFormats()
$ws.Range("E9").Value = '  -0.14%  '
$c = $ws.Range("D10")
$c.Value = "'" + '7.60'
$c.ClearFormats()
$ws.Range("E10").Value = '  +4.29%  '
$ws.Range("E11").Value = '  +2.60%  '
$c = $ws.Range("D12")
$c.Value = "'" + '0.387'
$c.ClearFormats()
$ws.Range("E12").Value = '  +2.95%  '
$ws.Range("D13").Value = '3.936.38'
$ws.Range("E13").Value = '  -0.06%  '
$ws.Range("E14").Value = '  +1.43%  '
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '3.361.92'
$ws.Range("E15").Value = '  -0.09%  '
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$c = $ws.Range("D16")
$c.Value = "'" + '0.0000171'
$c.ClearFormats()
$ws.Range("E16").Value = '  +0.58%  '
$c = $ws.Range("D17")
$c.Value = "'" + '24.96'
$c.ClearFormats()
$ws.Range("E17").Value = '  +2.51%  '
$ws.Range("D18").Value = '60.899.29'
$ws.Range("E18").Value = '  -2.94%  '
$c = $ws.Range("D19")
$c.Value = "'" + '13.85'
$c.ClearFormats()
$ws.Range("E19").Value = '  +6.06%  '
$c = $ws.Range("D20")
$c.Value = "'" + '9.37'
$c.ClearFormats()
$ws.Range("E20").Value = '  +1.62%  '
$c = $ws.Range("D21")
$c.Value = "'" + '5.73'
$c.ClearFormats()
$ws.Range("E21").Value = '  +2.13%  '
$c = $ws.Range("D22")
$c.Value = "'" + '372.05'
$c.ClearFormats()
$ws.Range("E22").Value = '  +0.25%  '
$c = $ws.Range("D23")
$c.Value = "'" + '0.571'
$c.ClearFormats()
$ws.Range("E23").Value = '  +3.46%  '
$ws.Range("D24").Value = '3.492.87'
$ws.Range("E24").Value = '  -0.44%  '
$ws.Range("E25").Value = '  -0.13%  '
$c = $ws.Range("D26")
$c.Value = "'" + '70.50'
$c.ClearFormats()
$ws.Range("E26").Value = '  -1.21%  '
$c = $ws.Range("D27")
$c.Value = "'" + '0.0000116'
$c.ClearFormats()
$ws.Range("E27").Value = '  +10.76%  '
$c = $ws.Range("D28")
$c.Value = "'" + '1.66'
$c.ClearFormats()
$ws.Range("E28").Value = '  +20.42%  '
$c = $ws.Range("D29")
$c.Value = "'" + '7.62'
$c.ClearFormats()
$ws.Range("E29").Value = '  +9.75%  '
$c = $ws.Range("D30")
$c.Value = "'" + '1.00'
$c.ClearFormats()
$ws.Range("E30").Value = '  -0.21%  '
$c = $ws.Range("D31")
$c.Value = "'" + '8.06'
$c.ClearFormats()
$ws.Range("E31").Value = '  +3.39%  '
$c = $ws.Range("D32")
$c.Value = "'" + '2.13'
$c.ClearFormats()
$ws.Range("E32").Value = '  +1.38%  '
$c = $ws.Range("D33")
$c.Value = "'" + '0.154'
$c.ClearFormats()
$ws.Range("E33").Value = '  +4.05%  '
$ws.Range("E34").Value = '  -0.06%  '
$ws.Range("D35").Value = '3.388.67'
$ws.Range("E35").Value = '  -0.48%  '
$c = $ws.Range("D36")
$c.Value = "'" + '23.26'
$c.ClearFormats()
$ws.Range("E36").Value = '  +2.74%  '
$c = $ws.Range("D37")
$c.Value = "'" + '5.52'
$c.ClearFormats()
$ws.Range("E37").Value = '  +2.27%  '
$c = $ws.Range("D38")
$c.Value = "'" + '6.89'
$c.ClearFormats()
$ws.Range("E38").Value = '  +4.23%  '
$c = $ws.Range("D39")
$c.Value = "'" + '1.54'
$c.ClearFormats()
$ws.Range("E39").Value = '  +4.94%  '
$c = $ws.Range("D40")
$c.Value = "'" + '162.19'
$c.ClearFormats()
$ws.Range("E40").Value = '  -2.40%  '
$c = $ws.Range("D41")
$c.Value = "'" + '0.0780'
$c.ClearFormats()
$ws.Range("E41").Value = '  +3.60%  '
$ws.Range("E42").Value = '  +0.17%  '
$ws.Range("B43").Value = 'OKB'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$c = $ws.Range("D43")
$c.Value = "'" + '41.22'
$c.ClearFormats()
$ws.Range("E43").Value = '  -0.77%  '
$ws.Range("B44").Value = 'Filecoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Range("D44")
$c.Value = "'" + '4.37'
$c.ClearFormats()
$ws.Range("E44").Value = '  +3.94%  '
$ws.Range("B45").Value = 'ONDO'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$c = $ws.Range("D45")
$c.Value = "'" + '1.20'
$c.ClearFormats()
$ws.Range("E45").Value = '  +10.35%  '
$c = $ws.Range("D46")
$c.Value = "'" + '0.755'
$c.ClearFormats()
$ws.Range("E46").Value = '  -1.09%  '
$c = $ws.Range("D47")
$c.Value = "'" + '1.59'
$c.ClearFormats()
$ws.Range("E47").Value = '  +3.32%  '
$ws.Range("B48").Value = 'Cosmos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c = $ws.Range("D48")
$c.Value = "'" + '6.94'
$c.ClearFormats()
$ws.Range("E48").Value = '  +5.63%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range("D49")
$c.Value = "'" + '22.82'
$c.ClearFormats()
$ws.Range("E49").Value = '  +1.49%  '
$c = $ws.Range("D50")
$c.Value = "'" + '23.23'
$c.ClearFormats()
$ws.Range("E50").Value = '  +15.01%  '
$c = $ws.Range("D51")
$c.Value = "'" + '2.40'
$c.ClearFormats()
$ws.Range("E51").Value = '  +13.20%  '
